# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions" (Sun Jan  7 09:43:32 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are plain text in the source sheet (e.g. "44.449.27",
# "0.0800"), but Excel's COM layer auto-coerces plain numeric-looking strings
# assigned through .Value into real numbers (dropping trailing zeros, introducing
# binary-float noise, etc.). Force the target range to Text first, write, then
# restore the default "Normal" style so no visible formatting changes stick.
function Set-TextValue {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "44.449.27"
$ws.Range("E2").Value = "  +0.96%  "

# Row 3
Set-TextValue $ws "D3" "2.236.55"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("E4").Value = "  +1.13%  "

# Row 5
Set-TextValue $ws "D5" "307.54"
$ws.Range("E5").Value = "  +0.83%  "

# Row 6
Set-TextValue $ws "D6" "93.60"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.19%  "

# Row 9
Set-TextValue $ws "D9" "0.520"
$ws.Range("E9").Value = "  -0.68%  "

# Row 10
Set-TextValue $ws "D10" "34.55"
$ws.Range("E10").Value = "  -1.07%  "

# Row 11
Set-TextValue $ws "D11" "0.0800"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
Set-TextValue $ws "D12" "7.18"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("E13").Value = "  +0.49%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws "D14" "0.830"
$ws.Range("E14").Value = "  +0.56%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D15" "2.211.32"
$ws.Range("E15").Value = "  -3.17%  "

# Row 16
Set-TextValue $ws "D16" "13.48"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
Set-TextValue $ws "D17" "44.096.30"
$ws.Range("E17").Value = "  +0.47%  "

# Row 18
Set-TextValue $ws "D18" "0.0₃0950"
$ws.Range("E18").Value = "  -0.69%  "

# Row 19
Set-TextValue $ws "D19" "6.32"
$ws.Range("E19").Value = "  +1.66%  "

# Row 20
Set-TextValue $ws "D20" "11.90"
$ws.Range("E20").Value = "  -3.05%  "

# Row 21
Set-TextValue $ws "D21" "65.63"
$ws.Range("E21").Value = "  +1.38%  "

# Row 22
Set-TextValue $ws "D22" "237.51"
$ws.Range("E22").Value = "  +0.60%  "

# Row 23
$ws.Range("E23").Value = "  +1.29%  "

# Row 24
Set-TextValue $ws "D24" "1.97"
$ws.Range("E24").Value = "  +0.99%  "

# Row 25
$ws.Range("E25").Value = "  -0.32%  "

# Row 26
$ws.Range("E26").Value = "  +4.03%  "

# Row 27
Set-TextValue $ws "D27" "9.75"
$ws.Range("E27").Value = "  -1.62%  "

# Row 28
Set-TextValue $ws "D28" "37.53"
$ws.Range("E28").Value = "  -0.90%  "

# Row 29
Set-TextValue $ws "D29" "5.89"
$ws.Range("E29").Value = "  -0.90%  "

# Row 30
Set-TextValue $ws "D30" "19.89"
$ws.Range("E30").Value = "  -0.56%  "

# Row 31
Set-TextValue $ws "D31" "153.67"
$ws.Range("E31").Value = "  -1.13%  "

# Row 32
Set-TextValue $ws "D32" "0.0793"
$ws.Range("E32").Value = "  -1.78%  "

# Row 33
$ws.Range("E33").Value = "  +0.25%  "

# Row 34
Set-TextValue $ws "D34" "3.09"
$ws.Range("E34").Value = "  -5.72%  "

# Row 35
Set-TextValue $ws "D35" "0.110"
$ws.Range("E35").Value = "  +1.73%  "

# Row 36
$ws.Range("E36").Value = "  +0.54%  "

# Row 37
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws "D38" "14.72"
$ws.Range("E38").Value = "  -3.88%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D39" "3.37"
$ws.Range("E39").Value = "  +0.54%  "

# Row 40
Set-TextValue $ws "D40" "3.75"
$ws.Range("E40").Value = "  -1.71%  "

# Row 41
$ws.Range("E41").Value = "  -0.82%  "

# Row 42
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
Set-TextValue $ws "D43" "1.769.95"
$ws.Range("E43").Value = "  +1.98%  "

# Row 44
$ws.Range("E44").Value = "  +1.79%  "

# Row 45
Set-TextValue $ws "D45" "78.78"
$ws.Range("E45").Value = "  -7.60%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D46" "98.40"
$ws.Range("E46").Value = "  -1.48%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws "D47" "4.88"
$ws.Range("E47").Value = "  -1.08%  "

# Row 48
Set-TextValue $ws "D48" "69.45"
$ws.Range("E48").Value = "  +0.64%  "

# Row 49
Set-TextValue $ws "D49" "8.07"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
Set-TextValue $ws "D50" "54.59"
$ws.Range("E50").Value = "  +0.69%  "

# Row 51
$ws.Range("E51").Value = "  +3.50%  "
